$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 1488.3334
$ws.Range("J20").Value = 3575
$ws.Range("L20").Value = 3575
$ws.Range("N20").Value = -4035
$ws.Range("H21").Value = 8016.6665
$ws.Range("I21").Value = 2025
$ws.Range("J21").Value = 20000
$ws.Range("K21").Value = 2025
$ws.Range("L21").Value = 20000
$ws.Range("M21").Value = -1557
$ws.Range("N21").Value = -20936
$ws.Range("H23").Value = 8016.6665
$ws.Range("I23").Value = 2025
$ws.Range("J23").Value = 20000
$ws.Range("K23").Value = 2025
$ws.Range("L23").Value = 20000
$ws.Range("M23").Value = -1791
$ws.Range("N23").Value = -20468
$ws.Range("H35").Value = 1488.3334
$ws.Range("J35").Value = 3575
$ws.Range("L35").Value = 3575
$ws.Range("N35").Value = -4333
$ws.Range("H43").Value = 2347.75
$ws.Range("I43").Value = 2130.3333
$ws.Range("J43").Value = 3000
$ws.Range("K43").Value = 2130.3333
$ws.Range("L43").Value = 3000
$ws.Range("M43").Value = -2061.3333
$ws.Range("N43").Value = -3138
$ws.Range("H64").Value = 2000
$ws.Range("J64").Value = 2000
$ws.Range("L64").Value = 2000
$ws.Range("N64").Value = -2496
$ws.Range("H67").Value = 2000
$ws.Range("J67").Value = 2000
$ws.Range("L67").Value = 2000
$ws.Range("N67").Value = -3716
$ws.Range("H92").Value = 772.5
$ws.Range("I92").Value = 1295
$ws.Range("K92").Value = 1295
$ws.Range("M92").Value = -47
$ws.Range("H100").Value = 1000
$ws.Range("I100").Value = 1000
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1000
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -459
$ws.Range("N100").ClearContents()
$ws.Range("H101").Value = 937.4286
$ws.Range("I101").Value = 595.6
$ws.Range("J101").Value = 1792
$ws.Range("K101").Value = 1786.8
$ws.Range("L101").Value = 5376
$ws.Range("M101").Value = -164.8000000000002
$ws.Range("N101").Value = -8620
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H129").Value = 3308.5454
$ws.Range("I129").Value = 1190.2858
$ws.Range("J129").Value = 4297.067
$ws.Range("K129").Value = 3570.8574
$ws.Range("L129").Value = 12891.201
$ws.Range("M129").Value = 1429.1426
$ws.Range("N129").Value = -22891.201
$ws.Range("H132").Value = 2814.6155
$ws.Range("I132").Value = 2550
$ws.Range("K132").Value = 7650
$ws.Range("M132").Value = -5120

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 331.2857
$ws.Range("I2").Value = 328.16666
$ws.Range("K2").Value = 328.16666
$ws.Range("M2").Value = -215.16666
$ws.Range("H4").Value = 2538
$ws.Range("J4").Value = 2547.5
$ws.Range("L4").Value = 2547.5
$ws.Range("N4").Value = -2779.5
$ws.Range("H74").Value = 3032.3333
$ws.Range("I74").Value = 3498.5
$ws.Range("K74").Value = 3498.5
$ws.Range("M74").Value = -2624.5
$ws.Range("H77").Value = 3032.3333
$ws.Range("I77").Value = 3498.5
$ws.Range("K77").Value = 17492.5
$ws.Range("M77").Value = -13124.5
$ws.Range("H97").Value = 838.4286
$ws.Range("I97").Value = 776
$ws.Range("K97").Value = 776
$ws.Range("M97").Value = -280
$ws.Range("H116").Value = 331.2857
$ws.Range("I116").Value = 328.16666
$ws.Range("K116").Value = 328.16666
$ws.Range("M116").Value = 1965.83334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 331.2857
$ws.Range("I3").Value = 328.16666
$ws.Range("K3").Value = 328.16666
$ws.Range("M3").Value = -214.16666
$ws.Range("H86").Value = 3239.9
$ws.Range("I86").Value = 3112.5
$ws.Range("J86").Value = 3749.5
$ws.Range("K86").Value = 3112.5
$ws.Range("L86").Value = 3749.5
$ws.Range("M86").Value = -1989.5
$ws.Range("N86").Value = -5995.5
$ws.Range("H89").Value = 3239.9
$ws.Range("I89").Value = 3112.5
$ws.Range("J89").Value = 3749.5
$ws.Range("K89").Value = 15562.5
$ws.Range("L89").Value = 18747.5
$ws.Range("M89").Value = -9946.5
$ws.Range("N89").Value = -29979.5
$ws.Range("H94").Value = 1124.2593
$ws.Range("I94").Value = 1163.0435
$ws.Range("J94").Value = 901.25
$ws.Range("K94").Value = 1163.0435
$ws.Range("L94").Value = 901.25
$ws.Range("M94").Value = -712.0435
$ws.Range("N94").Value = -1803.25
$ws.Range("H105").Value = 2957.8
$ws.Range("J105").Value = 3499.5
$ws.Range("L105").Value = 3499.5
$ws.Range("N105").Value = -6993.5
$ws.Range("H134").Value = 4336.2144
$ws.Range("I134").Value = 4448.56
$ws.Range("K134").Value = 13345.68
$ws.Range("M134").Value = -10810.68

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 496.33334
$ws.Range("I22").Value = 199.5
$ws.Range("K22").Value = 199.5
$ws.Range("M22").Value = 150.5
$ws.Range("H31").Value = 4623.9375
$ws.Range("I31").Value = 4256.6
$ws.Range("K31").Value = 4256.6
$ws.Range("M31").Value = -3961.6
$ws.Range("H34").Value = 4623.9375
$ws.Range("I34").Value = 4256.6
$ws.Range("K34").Value = 4256.6
$ws.Range("M34").Value = -4054.6
$ws.Range("H56").Value = 4850
$ws.Range("I56").Value = 4700
$ws.Range("J56").Value = 5000
$ws.Range("K56").Value = 4700
$ws.Range("L56").Value = 5000
$ws.Range("M56").Value = -3855
$ws.Range("N56").Value = -6690
$ws.Range("H134").Value = 3888.6365
$ws.Range("I134").Value = 2597.375
$ws.Range("K134").Value = 7792.125
$ws.Range("M134").Value = -5257.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 64.71429
$ws.Range("I2").Value = 88.416664
$ws.Range("J2").Value = 33.11111
$ws.Range("K2").Value = 530.499984
$ws.Range("L2").Value = 198.66666
$ws.Range("M2").Value = -417.499984
$ws.Range("N2").Value = -424.66666
$ws.Range("H140").Value = 1945.375
$ws.Range("I140").Value = 993.8333
$ws.Range("J140").Value = 4800
$ws.Range("K140").Value = 2981.4999
$ws.Range("L140").Value = 14400
$ws.Range("M140").Value = 2198.5001
$ws.Range("N140").Value = -24760

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 138
$ws.Range("I107").Value = 138
$ws.Range("K107").Value = 138
$ws.Range("M107").Value = 1782

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 783.8
$ws.Range("I22").Value = 771.2857
$ws.Range("J22").Value = 813
$ws.Range("K22").Value = 771.2857
$ws.Range("L22").Value = 813
$ws.Range("M22").Value = -476.2857
$ws.Range("N22").Value = -1403
$ws.Range("H27").Value = 783.8
$ws.Range("I27").Value = 771.2857
$ws.Range("J27").Value = 813
$ws.Range("K27").Value = 771.2857
$ws.Range("L27").Value = 813
$ws.Range("M27").Value = -664.2857
$ws.Range("N27").Value = -1027
$ws.Range("H95").Value = 16750
$ws.Range("J95").Value = 16750
$ws.Range("L95").Value = 16750
$ws.Range("N95").Value = -22242
$ws.Range("H132").Value = 15698.866
$ws.Range("I132").Value = 16727.818
$ws.Range("J132").Value = 12869.25
$ws.Range("K132").Value = 50183.454
$ws.Range("L132").Value = 38607.75
$ws.Range("M132").Value = -47653.454
$ws.Range("N132").Value = -43667.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 16737.5
$ws.Range("I58").Value = 7316.6665
$ws.Range("J58").Value = 45000
$ws.Range("K58").Value = 7316.6665
$ws.Range("L58").Value = 45000
$ws.Range("M58").Value = -7008.6665
$ws.Range("N58").Value = -45616
$ws.Range("H81").Value = 368.26666
$ws.Range("I81").Value = 371.07693
$ws.Range("K81").Value = 742.15386
$ws.Range("M81").Value = 318.84614
$ws.Range("H84").Value = 368.26666
$ws.Range("I84").Value = 371.07693
$ws.Range("K84").Value = 3710.7693
$ws.Range("M84").Value = 1593.2307
$ws.Range("H132").Value = 896.4286
$ws.Range("I132").Value = 643.75
$ws.Range("J132").Value = 1233.3334
$ws.Range("K132").Value = 1931.25
$ws.Range("L132").Value = 3700.0002
$ws.Range("M132").Value = 598.75
$ws.Range("N132").Value = -8760.0002
$ws.Range("H135").Value = 52570
$ws.Range("J135").Value = 58331.668
$ws.Range("L135").Value = 58331.668
$ws.Range("N135").Value = -68471.668
$ws.Range("H136").Value = 1715.7391
$ws.Range("I136").Value = 1657.3636
$ws.Range("K136").Value = 4972.0908
$ws.Range("M136").Value = -2422.0908
$ws.Range("H137").Value = 39000
$ws.Range("I137").Value = 39000
$ws.Range("K137").Value = 39000
$ws.Range("M137").Value = -33900
